# "cleaned defensive actions data"
#
# The sheet was originally written by pandas as a two-row MultiIndex header
# (row 1 = top level, with merged cells + "Unnamed: n_level_0" placeholders;
# row 2 = the real per-column labels). This cleans it up the way Excel does
# when a user flattens that header by hand:
#   * unmerge the three merged header spans on row 1
#   * replace row 1 with the real, flattened column headers (incl. two new
#     ones - "Player ID" for col A and a short "Cha" for the challenges
#     group - and "90s" instead of "Min")
#   * hide the old (now redundant) per-column header row (row 2)
#   * hide the always-blank spacer row (row 3)
#   * hide the summary/total row at the bottom (row 20)
#   * make the previously-omitted "Tkl%" column explicit (0) on the rows
#     where pandas skipped writing a literal zero
#   * leave the final selection on O21, like the saved workbook shows

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 used to be three merged spans (H1:L1, M1:P1, Q1:S1) holding a single
# group label each. Unmerge them so every column in row 1 can carry its own
# real header text.
$ws.Range("H1:L1").UnMerge()
$ws.Range("M1:P1").UnMerge()
$ws.Range("Q1:S1").UnMerge()

# Flattened header text for row 1, column A (1) through W (23).
$headers = @(
    "Player ID", "Player", "#", "Nation", "Pos", "Age", "90s",
    "Tkl", "TklW", "Def 3rd", "Mid 3rd", "Att 3rd",
    "Cha", "Att", "Tkl%", "Lost",
    "Blocks", "Sh", "Pass",
    "Int", "Tkl+Int", "Clr", "Err"
)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Row 2 (the old sub-header row) and row 3 (a blank spacer row) are no
# longer needed for display, so hide them instead of deleting them.
$ws.Rows.Item(2).Hidden = $true
$ws.Rows.Item(3).Hidden = $true

# A handful of rows were missing an explicit 0 in the "Tkl%" column (O) -
# fill those in.
$tklPctRows = @(5, 8, 11, 13, 19)
foreach ($r in $tklPctRows) {
    $ws.Cells.Item($r, 15).Value = 0
}

# The trailing "16 Players" totals row is hidden in the cleaned sheet.
$ws.Rows.Item(20).Hidden = $true

# Leave the selection where the saved workbook shows it.
$ws.Range("O21").Select()
